# Refresh the cryptos price/volume table with the latest scraped figures.
# (Values are stored as plain text in the sheet, so numeric-looking prices
# are written with a leading quote to stop Excel from re-typing them as
# numbers, then the cell style is reset to Normal so no stray number format
# is left behind.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '34.680.59'
$ws.Range('E2').Value = '  -2.16%  '
$ws.Range('D3').Value = '1.810.75'
$ws.Range('E3').Value = '  -1.74%  '
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').Value = "'232.33"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('E6').Value = '  -0.96%  '
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('D8').Value = "'39.32"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -9.57%  '
$ws.Range('E9').Value = '  +5.61%  '
$ws.Range('D10').Value = "'0.0681"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.50%  '
$ws.Range('E11').Value = '  -2.10%  '
$ws.Range('D12').Value = '2.072.04'
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('D13').Value = "'0.673"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('D14').Value = '1.816.23'
$ws.Range('D15').Value = "'11.20"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.10%  '
$ws.Range('E16').Value = '  -2.24%  '
$ws.Range('D17').Value = '34.679.80'
$ws.Range('E17').Value = '  -2.02%  '
$ws.Range('D18').Value = "'69.42"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('E19').Value = '  -1.57%  '
$ws.Range('D20').Value = "'240.36"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.85%  '
$ws.Range('D21').Value = "'11.92"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').Value = "'4.70"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('E24').Value = '  +1.89%  '
$ws.Range('D25').Value = "'171.99"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('E26').Value = '  -2.75%  '
$ws.Range('D27').Value = "'17.22"
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Value = "'0.120"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.94%  '
$ws.Range('E29').Value = '  +2.20%  '
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('D31').Value = "'4.02"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.35%  '
$ws.Range('E32').Value = '  -0.22%  '
$ws.Range('E33').Value = '  -2.73%  '
$ws.Range('D34').Value = "'1.29"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +18.61%  '
$ws.Range('D35').Value = "'1.79"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.84%  '
$ws.Range('E36').Value = '  +1.50%  '
$ws.Range('D37').Value = "'91.52"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.91%  '
$ws.Range('E38').Value = '  +4.33%  '
$ws.Range('D39').Value = '1.325.99'
$ws.Range('E39').Value = '  -1.54%  '
$ws.Range('E40').Value = '  -1.02%  '
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('D42').Value = "'0.963"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.57%  '
$ws.Range('D43').Value = "'14.14"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.69%  '
$ws.Range('E44').Value = '  -8.89%  '
$ws.Range('E45').Value = '  -4.85%  '
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('D47').Value = "'0.0513"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.04%  '
$ws.Range('D48').Value = '1.999.98'
$ws.Range('E48').Value = '  -0.44%  '
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('E50').Value = '  +7.08%  '
$ws.Range('D51').Value = "'98.35"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.58%  '
